$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$sh = $s.Shapes.Item("Rectangle 28")
$sh.Left = 633.4003149606299
$sh.Top = 158.10031896062992
$sh.Width = 83.13007874015749
$sh.Height = 102.99937057874016

$sh = $s.Shapes.Item("Rectangle 27")
$sh.Left = 546.8270866141733
$sh.Top = 158.10031896062992
$sh.Width = 83.98669291338582
$sh.Height = 102.99811173622047

$sh = $s.Shapes.Item("Rectangle 22")
$sh.Left = 662.0591338582677
$sh.Top = 394.29748031496064

$sh = $s.Shapes.Item("TextBox 6")
$sh.Left = 117.22464566929133
$sh.Top = 363.3995275590551

$sh = $s.Shapes.Item("TextBox 7")
$sh.Left = 508.7928346456693
$sh.Top = 364.09898377795275

$sh = $s.Shapes.Item("Straight Connector 11")
$sh.Left = 267.8641815283464
$sh.Top = 383.9981232362205

$sh = $s.Shapes.Item("Straight Connector 14")
$sh.Left = 645.8046570692912
$sh.Top = 384.6981964763779

$sh = $s.Shapes.Item("Chord 17")
$sh.Left = 297.68661417322835
$sh.Top = 270.6994488188976

$sh = $s.Shapes.Item("Oval 18")
$sh.Left = 299.43653873307085
$sh.Top = 393.5985039370079

$sh = $s.Shapes.Item("Chord 19")
$sh.Left = 298.07386786771656
$sh.Top = 393.5985039370079

$sh = $s.Shapes.Item("Chord 20")
$sh.Left = 697.2122047244095
$sh.Top = 270.6994488188976

$sh = $s.Shapes.Item("Chord 21")
$sh.Left = 697.2122047244095
$sh.Top = 394.298031496063

$sh = $s.Shapes.Item("Oval 25")
$sh.Left = 579.9773255346456
$sh.Top = 158.10094488188977
$sh.Width = 105.36645889291339
$sh.Height = 102.99937057874016

$sh = $s.Shapes.Item("Chord 26")
$sh.Left = 578.3966369732284
$sh.Top = 158.10031896062992
$sh.Width = 108.74968503937008
$sh.Height = 102.99937057874016

$sh = $s.Shapes.Item("TextBox 1")
$sh.Left = 297.68661417322835
$sh.Top = 189.00016028031496
$sh.Width = 235.98417322834646
$sh.Height = 41.198425196850394
